$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue $ws "D2" "59.907.13"
Set-TextValue $ws "E2" "  +0.09%  "
Set-TextValue $ws "D3" "2.377.68"
Set-TextValue $ws "E3" "  -1.39%  "
Set-TextValue $ws "E4" "  +0.09%  "
Set-TextValue $ws "D5" "556.62"
Set-TextValue $ws "E5" "  +0.97%  "
Set-TextValue $ws "D6" "133.67"
Set-TextValue $ws "E6" "  -2.47%  "
Set-TextValue $ws "E7" "  +0.08%  "
Set-TextValue $ws "D8" "0.585"
Set-TextValue $ws "E8" "  -0.96%  "
Set-TextValue $ws "E9" "  +0.79%  "
Set-TextValue $ws "D10" "5.65"
Set-TextValue $ws "E10" "  -0.35%  "
Set-TextValue $ws "E11" "  +1.05%  "
Set-TextValue $ws "D12" "0.343"
Set-TextValue $ws "E12" "  -3.16%  "
Set-TextValue $ws "D13" "24.41"
Set-TextValue $ws "E13" "  -3.70%  "
Set-TextValue $ws "D14" "2.805.58"
Set-TextValue $ws "E14" "  -1.16%  "
Set-TextValue $ws "D15" "59.932.22"
Set-TextValue $ws "E15" "  +0.27%  "
Set-TextValue $ws "E16" "  +0.48%  "
Set-TextValue $ws "D17" "2.384.53"
Set-TextValue $ws "E17" "  -1.63%  "
Set-TextValue $ws "D18" "11.12"
Set-TextValue $ws "E18" "  -1.99%  "
Set-TextValue $ws "D19" "4.49"
Set-TextValue $ws "E19" "  +1.64%  "
Set-TextValue $ws "D20" "321.42"
Set-TextValue $ws "E20" "  -2.16%  "
Set-TextValue $ws "D21" "6.66"
Set-TextValue $ws "E21" "  -0.08%  "
Set-TextValue $ws "D22" "0.999"
Set-TextValue $ws "E22" "  +0.01%  "
Set-TextValue $ws "D23" "64.15"
Set-TextValue $ws "E23" "  -3.38%  "
Set-TextValue $ws "E24" "  +0.40%  "
Set-TextValue $ws "E25" "  -0.02%  "
Set-TextValue $ws "D26" "8.44"
Set-TextValue $ws "E26" "  -2.60%  "
Set-TextValue $ws "D27" "1.37"
Set-TextValue $ws "E27" "  +0.30%  "
Set-TextValue $ws "E28" "  +2.05%  "
Set-TextValue $ws "D29" "0.0₃0760"
Set-TextValue $ws "E29" "  -1.54%  "
Set-TextValue $ws "D30" "169.84"
Set-TextValue $ws "E30" "  +0.56%  "
Set-TextValue $ws "D31" "6.07"
Set-TextValue $ws "E31" "  +0.80%  "
Set-TextValue $ws "E32" "  +10.84%  "
Set-TextValue $ws "D33" "0.401"
Set-TextValue $ws "E33" "  -1.24%  "
Set-TextValue $ws "D34" "18.15"
Set-TextValue $ws "E34" "  -2.48%  "
Set-TextValue $ws "B35" "USDe"
Set-TextValue $ws "C35" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D35" "0.999"
Set-TextValue $ws "E35" "  +0.00%  "
Set-TextValue $ws "B36" "ImmutableX"
Set-TextValue $ws "C36" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D36" "1.33"
Set-TextValue $ws "E36" "  +1.47%  "
Set-TextValue $ws "E37" "  +0.14%  "
Set-TextValue $ws "D38" "4.14"
Set-TextValue $ws "E38" "  -1.25%  "
Set-TextValue $ws "E39" "  -0.86%  "
Set-TextValue $ws "D40" "318.35"
Set-TextValue $ws "E40" "  -0.45%  "
Set-TextValue $ws "D41" "38.61"
Set-TextValue $ws "E41" "  -1.97%  "
Set-TextValue $ws "D42" "144.91"
Set-TextValue $ws "E42" "  +3.64%  "
Set-TextValue $ws "D43" "3.53"
Set-TextValue $ws "E43" "  -3.37%  "
Set-TextValue $ws "E44" "  +0.41%  "
Set-TextValue $ws "D45" "19.74"
Set-TextValue $ws "E45" "  +1.07%  "
Set-TextValue $ws "D46" "0.0512"
Set-TextValue $ws "E46" "  -0.59%  "
Set-TextValue $ws "D47" "0.569"
Set-TextValue $ws "E47" "  -1.41%  "
Set-TextValue $ws "D48" "0.0218"
Set-TextValue $ws "E48" "  -2.18%  "
Set-TextValue $ws "E49" "  +0.11%  "
Set-TextValue $ws "B50" "ZEEBU"
Set-TextValue $ws "C50" "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
Set-TextValue $ws "D50" "4.67"
Set-TextValue $ws "E50" "  -0.06%  "
Set-TextValue $ws "B51" "dogwifhat"
Set-TextValue $ws "C51" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D51" "1.54"
Set-TextValue $ws "E51" "  -1.70%  "
